# adding averages and more checks
$wb = $excel.ActiveWorkbook

$trainingWs = $wb.Worksheets.Item("Training Dashboard")
$examWs     = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------------
# Header rows: make the header text bold white (on top of the existing blue
# fill / border / centered alignment) on both dashboards.
# ---------------------------------------------------------------------------
$trainingWs.Range("A2:K2").Font.Color = 16777215
$examWs.Range("A2:G2").Font.Color = 16777215

# ---------------------------------------------------------------------------
# Training Dashboard: recalculate "PERIOD TO EXPIRE" (column H) now that the
# "LAST UPDATE" (column I) moved forward from 08-Sep-2025 to 16-Sep-2025,
# i.e. 8 fewer days remain until expiry.
# ---------------------------------------------------------------------------
for ($row = 3; $row -le 23; $row++) {
    $periodCell = $trainingWs.Cells.Item($row, 8)
    $periodCell.Value = $periodCell.Value() - 8

    $lastUpdateCell = $trainingWs.Cells.Item($row, 9)
    $lastUpdateCell.Value = "'16-Sep-2025"
}

# ---------------------------------------------------------------------------
# Exam Dashboard: widen the COMMENTS column and refresh its remark text.
# ---------------------------------------------------------------------------
$examWs.Columns.Item(5).ColumnWidth = 14.17
$examWs.Range("E3").Value = "date is valid"
